# Update column G ("K") values on the active sheet (data rows 2-11)
# to reflect the regenerated save_data (K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 5
    4  = 6
    5  = 16
    6  = 2
    7  = 3
    8  = 3
    9  = 2
    10 = 2
    11 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
